$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Strip the "[1]"/"[2]" citation markers from the organism names (col B,
#    rows 4-15) now that the References section is gone.
# ---------------------------------------------------------------------------
$orgNames = @{
  4  = "Gordonia phage GAL1"
  5  = "WS1 bacterium JGI 0000059-K21"
  6  = "Astrammina rara"
  7  = "Nosema ceranae"
  8  = "Cryptosporidium parvum Iowa II"
  9  = "Spironucleus salmonicida"
  10 = "Tieghemostelium lacteum"
  11 = "Fusarium graminearum PH-1"
  12 = "Salpingoeca rosetta"
  13 = "Chondrus crispus"
  14 = "Kappaphycus alvarezii"
  15 = "Strongylocentrotus purpuratus"
}
foreach ($r in $orgNames.Keys) {
  $ws.Cells.Item($r, 2).Value = $orgNames[$r]
}

# ---------------------------------------------------------------------------
# 2) Title: "Datasets Description" -> "DNA Datasets Description"
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "DNA Datasets Description"

# ---------------------------------------------------------------------------
# 3) Column widths: A & B resized, new widths for H & I
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 21.42578125
$ws.Columns.Item(2).ColumnWidth = 31.42578125
$ws.Columns.Item(8).ColumnWidth = 10.5703125
$ws.Columns.Item(9).ColumnWidth = 11.28515625

# ---------------------------------------------------------------------------
# 4) Remove the old "References" block (rows 18-22) entirely, then rebuild
#    a new "RNA Datasets Description" block underneath the DNA table.
# ---------------------------------------------------------------------------
$ws.Range("A18:I22").EntireRow.Delete()

# Row 16 stays a blank spacer row (already blank after delete - nothing to do)

# Row 17: section title, merged A17:I17, same look as the "DNA Datasets
# Description" style used for the table header (bold Times New Roman 12,
# centered, bottom border).
$ws.Range("A2:D2").Copy()
$ws.Range("A17:I17").PasteSpecial(-4122)
$ws.Range("A17").Value = "RNA Datasets Description"
$ws.Range("A17:I17").Borders.Item(9).LineStyle = 1
$ws.Range("A17:I17").Borders.Item(9).Weight = 2
$ws.Range("A17:I17").Borders.Item(8).LineStyle = -4142
$ws.Range("A17:I17").Merge()

# Row 18/19: header row (merged vertically for A, C:H, I ; B is single row)
$ws.Range("A2:D2").Copy()
$ws.Range("A18:I19").PasteSpecial(-4122)

$ws.Range("A18").Value = "Datasets"
$ws.Range("B18").Value = "O.Size (B)"
$ws.Range("C18").Value = "Source"
$ws.Range("I18").Value = "Accessed Date"

$ws.Range("A18:A19").WrapText = $true
$ws.Range("C18:H19").WrapText = $true
$ws.Range("I18:I19").WrapText = $true

# "Accessed Date" uses the workbook's bold default font (Calibri 11) rather
# than the Times New Roman used elsewhere in the header.
$ws.Range("I18:I19").Font.Name = "Calibri"
$ws.Range("I18:I19").Font.Size = 11
$ws.Range("I18:I19").Font.Bold = $true

$ws.Range("A18:I18").Borders.Item(8).LineStyle = 1
$ws.Range("A18:I18").Borders.Item(8).Weight = 2
$ws.Range("A19:I19").Borders.Item(9).LineStyle = 1
$ws.Range("A19:I19").Borders.Item(9).Weight = 2

$ws.Range("A18:A19").Merge()
$ws.Range("B18:B19").Merge()
$ws.Range("C18:H19").Merge()
$ws.Range("I18:I19").Merge()

$ws.Rows.Item(18).RowHeight = 15
$ws.Rows.Item(19).RowHeight = 15.75

# Row 20: first data row (SILVA 132 LSURef)
$ws.Range("A4").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = "SILVA 132 LSURef "
$ws.Range("A20").Borders.Item(8).LineStyle = 1
$ws.Range("A20").Borders.Item(8).Weight = 2

$ws.Range("B20").NumberFormat = "#,##0"
$ws.Range("B20").HorizontalAlignment = -4108
$ws.Range("B20").VerticalAlignment = -4108
$ws.Range("B20").Value = 610296406

$ws.Range("I20").Font.Name = "Times New Roman"
$ws.Range("I20").Font.Size = 12
$ws.Range("I20").Font.Bold = $false
$ws.Range("I20").HorizontalAlignment = -4131
$ws.Range("I20").VerticalAlignment = -4108
$ws.Range("I20").Value = "19/08/2022"

$ws.Range("C20:H20").Merge()
$ws.Hyperlinks.Add($ws.Range("C20"), "https://ftp.arb-silva.de/release%20132/Exports/SILVA_132_LSURef_tax_silva.fasta.gz", "", "", "https://ftp.arb-silva.de/release 132/Exports/SILVA 132 LSURef tax silva.fasta.gz")
$ws.Range("C20:H20").WrapText = $true
$ws.Range("C20:H20").HorizontalAlignment = -4131
$ws.Range("C20:H20").VerticalAlignment = -4108

$ws.Rows.Item(20).RowHeight = 15

# Row 21: second data row (SILVA 132 SSURef Nr99)
$ws.Range("A4").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value = "SILVA 132 SSURef Nr99"
$ws.Range("A21").Borders.Item(9).LineStyle = 1
$ws.Range("A21").Borders.Item(9).Weight = 2

$ws.Range("B21").NumberFormat = "#,##0"
$ws.Range("B21").HorizontalAlignment = -4108
$ws.Range("B21").VerticalAlignment = -4108
$ws.Range("B21").Borders.Item(9).LineStyle = 1
$ws.Range("B21").Borders.Item(9).Weight = 2
$ws.Range("B21").Value = 1108994702

$ws.Range("I21").Font.Name = "Times New Roman"
$ws.Range("I21").Font.Size = 12
$ws.Range("I21").Font.Bold = $false
$ws.Range("I21").HorizontalAlignment = -4131
$ws.Range("I21").VerticalAlignment = -4108
$ws.Range("I21").Borders.Item(9).LineStyle = 1
$ws.Range("I21").Borders.Item(9).Weight = 2
$ws.Range("I21").Value = "19/08/2022"

$ws.Range("C21:H21").Merge()
$ws.Hyperlinks.Add($ws.Range("C21"), "https://ftp.arb-silva.de/release%20132/Exports/SILVA_132_SSURef_Nr99_tax_silva.fasta.gz", "", "", "https://ftp.arb-silva.de/release 132/Exports/SILVA 132 SSURef Nr99 tax silva.fasta.gz")
$ws.Range("C21:H21").WrapText = $true
$ws.Range("C21:H21").HorizontalAlignment = -4131
$ws.Range("C21:H21").VerticalAlignment = -4108
$ws.Range("C21:H21").Borders.Item(9).LineStyle = 1
$ws.Range("C21:H21").Borders.Item(9).Weight = 2

$ws.Rows.Item(21).RowHeight = 15

# Rows 22-23: trailing blank rows (Times New Roman 12, wrap, vertical center)
$ws.Range("A22:G23").Font.Name = "Times New Roman"
$ws.Range("A22:G23").Font.Size = 12
$ws.Range("A22:G23").WrapText = $true
$ws.Range("A22:G23").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 5) Sheet view: drop the frozen/scrolled topLeftCell and move the selection.
# ---------------------------------------------------------------------------
$ws.Range("H11").Select()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
